# "added signin and async wait"
# Update the SignIn test-data sheet with a new login (username/password)
# and leave that sheet as the active/selected one (it was previously
# "Profile" that was active).

$wb = $excel.ActiveWorkbook

$signIn = $wb.Worksheets.Item("SignIn")

# New credentials for the SignIn test data (note trailing spaces, matching
# the values used by the updated automation test).
$signIn.Range("B2").Value = "mvpstudio.qa@gmail.com "
$signIn.Range("C2").Value = "SydneyQa2018 "

# Make the SignIn sheet the active tab/selection.
$signIn.Activate()
